# "Added Week 15 simulations"
#
# 1) Rushing sheet: update a handful of cumulative simulation stats.
# 2) Receiving sheet: update a handful of cumulative simulation stats,
#    insert a new row for E.Winston at row 17 (pushing the players that
#    were on rows 17-20 down to rows 18-21), and set that new row's data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Rushing
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2 - T.Hill
$rushing.Cells.Item(2, 3).Value = 10
$rushing.Cells.Item(2, 4).Value = 17
$rushing.Cells.Item(2, 5).Value = 15
$rushing.Cells.Item(2, 6).Value = 14

# Row 4 - A.Kamara
$rushing.Cells.Item(4, 3).Value = 98
$rushing.Cells.Item(4, 4).Value = 61
$rushing.Cells.Item(4, 5).Value = 14
$rushing.Cells.Item(4, 6).Value = 29

# Row 5 - T.Jones
$rushing.Cells.Item(5, 3).Value = 22
$rushing.Cells.Item(5, 4).Value = 11
$rushing.Cells.Item(5, 5).Value = 6
$rushing.Cells.Item(5, 6).Value = 4

# ---------------------------------------------------------------------
# Receiving
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2 - A.Kamara
$receiving.Cells.Item(2, 3).Value = 44
$receiving.Cells.Item(2, 4).Value = 33
$receiving.Cells.Item(2, 5).Value = 5
$receiving.Cells.Item(2, 6).Value = 3
$receiving.Cells.Item(2, 7).Value = 10
$receiving.Cells.Item(2, 8).Value = 8

# Row 3 - T.Jones
$receiving.Cells.Item(3, 3).Value = 7
$receiving.Cells.Item(3, 4).Value = 5
$receiving.Cells.Item(3, 5).Value = 0
$receiving.Cells.Item(3, 6).Value = 0
$receiving.Cells.Item(3, 7).Value = 1
$receiving.Cells.Item(3, 8).Value = 1

# Row 9 - T.Smith
$receiving.Cells.Item(9, 3).Value = 29
$receiving.Cells.Item(9, 4).Value = 20
$receiving.Cells.Item(9, 5).Value = 12
$receiving.Cells.Item(9, 6).Value = 5
$receiving.Cells.Item(9, 7).Value = 5
$receiving.Cells.Item(9, 8).Value = 3

# Row 10 - M.Callaway
$receiving.Cells.Item(10, 3).Value = 38
$receiving.Cells.Item(10, 4).Value = 24
$receiving.Cells.Item(10, 5).Value = 20
$receiving.Cells.Item(10, 6).Value = 5
$receiving.Cells.Item(10, 7).Value = 7
$receiving.Cells.Item(10, 8).Value = 5

# New player E.Winston joins the roster at row 17, pushing A.Trautman,
# J.Johnson, G.Griffin and N.Vannett each down by one row.
$receiving.Rows.Item(17).Insert()

# Match the bold/centered/bordered style used by the rest of column A.
$newA17 = $receiving.Cells.Item(17, 1)
$newA17.Font.Bold = $true
$newA17.HorizontalAlignment = -4108
$newA17.VerticalAlignment = -4160
$newA17.Borders.LineStyle = 1

$receiving.Cells.Item(17, 1).Value = 15
$receiving.Cells.Item(17, 2).Value = "E.Winston"
$receiving.Cells.Item(17, 3).Value = 1
$receiving.Cells.Item(17, 4).Value = 1
$receiving.Cells.Item(17, 5).Value = 0
$receiving.Cells.Item(17, 6).Value = 0
$receiving.Cells.Item(17, 7).Value = 0
$receiving.Cells.Item(17, 8).Value = 0

# Row 18 - A.Trautman (formerly row 17)
$receiving.Cells.Item(18, 1).Value = 16
$receiving.Cells.Item(18, 3).Value = 36
$receiving.Cells.Item(18, 4).Value = 22
$receiving.Cells.Item(18, 5).Value = 2
$receiving.Cells.Item(18, 6).Value = 2
$receiving.Cells.Item(18, 7).Value = 7
$receiving.Cells.Item(18, 8).Value = 3

# Row 19 - J.Johnson (formerly row 18)
$receiving.Cells.Item(19, 1).Value = 17
$receiving.Cells.Item(19, 3).Value = 12
$receiving.Cells.Item(19, 4).Value = 8
$receiving.Cells.Item(19, 5).Value = 7
$receiving.Cells.Item(19, 6).Value = 3
$receiving.Cells.Item(19, 7).Value = 4
$receiving.Cells.Item(19, 8).Value = 3

# Row 20 - G.Griffin (formerly row 19)
$receiving.Cells.Item(20, 1).Value = 18
$receiving.Cells.Item(20, 3).Value = 4
$receiving.Cells.Item(20, 4).Value = 4
$receiving.Cells.Item(20, 5).Value = 0
$receiving.Cells.Item(20, 6).Value = 0
$receiving.Cells.Item(20, 7).Value = 0
$receiving.Cells.Item(20, 8).Value = 0

# Row 21 - N.Vannett (formerly row 20)
$receiving.Cells.Item(21, 1).Value = 19
$receiving.Cells.Item(21, 2).Value = "N.Vannett"
$receiving.Cells.Item(21, 3).Value = 10
$receiving.Cells.Item(21, 4).Value = 6
$receiving.Cells.Item(21, 5).Value = 2
$receiving.Cells.Item(21, 6).Value = 2
$receiving.Cells.Item(21, 7).Value = 3
$receiving.Cells.Item(21, 8).Value = 1
